$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 22:43"

# Row 4: 'Estados Unidos'
$ws.Range("B4").Value = 4846626
$ws.Range("C4").Value = 32979
$ws.Range("D4").Value = 2403683
$ws.Range("E4").Value = 2284237
$ws.Range("G4").Value = 341
$ws.Range("H4").Value = 158706

# Row 8: 'Sudafrica'
$ws.Range("B8").Value = 516862
$ws.Range("C8").Value = 5377
$ws.Range("D8").Value = 358037
$ws.Range("E8").Value = 150286
$ws.Range("G8").Value = 173
$ws.Range("H8").Value = 8539

# Row 21: 'Alemania'
$ws.Range("B21").Value = 212264
$ws.Range("C21").Value = 802
$ws.Range("E21").Value = 9432
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = 9232

# Row 36: 'Israel'
$ws.Range("B36").Value = 74430
$ws.Range("C36").Value = 1615
$ws.Range("D36").Value = 47571
$ws.Range("E36").Value = 26313

# Row 52: 'Barein'
$ws.Range("E52").Value = 2720
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 150

# Row 70: 'Etiopia' -> 'Costa Rica'
$ws.Range("A70").Value = 'Costa Rica'
$ws.Range("B70").Value = 19402
$ws.Range("C70").Value = 427
$ws.Range("D70").Value = 4689
$ws.Range("E70").Value = 14542
$ws.Range("G70").Value = 9
$ws.Range("H70").Value = 171

# Row 71: 'Costa Rica' -> 'Etiopia'
$ws.Range("A71").Value = 'Etiopia'
$ws.Range("B71").Value = 19289
$ws.Range("C71").Value = 583
$ws.Range("D71").Value = 7931
$ws.Range("E71").Value = 11022
$ws.Range("G71").Value = 26
$ws.Range("H71").Value = 336

# Row 76: 'Costa de Marfil'
$ws.Range("B76").Value = 16220
$ws.Range("C76").Value = 38
$ws.Range("D76").Value = 11887
$ws.Range("E76").Value = 4231

# Row 90: 'Tayikistan' -> 'Gabon'
$ws.Range("A90").Value = 'Gabon'
$ws.Range("B90").Value = 7646
$ws.Range("C90").Value = 115
$ws.Range("D90").Value = 5408
$ws.Range("E90").Value = 2187
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 51

# Row 91: 'Gabon' -> 'Tayikistan'
$ws.Range("A91").Value = 'Tayikistan'
$ws.Range("B91").Value = 7538
$ws.Range("C91").Value = 43
$ws.Range("D91").Value = 6317
$ws.Range("E91").Value = 1160
$ws.Range("H91").Value = 61

# Row 97: 'Mauritania'
$ws.Range("B97").Value = 6382
$ws.Range("C97").Value = 59
$ws.Range("D97").Value = 5174

# Row 108: 'Malaui'
$ws.Range("B108").Value = 4272
$ws.Range("C108").Value = 41
$ws.Range("D108").Value = 1945
$ws.Range("E108").Value = 2204

# Row 113: 'Tailandia' -> 'Congo'
$ws.Range("A113").Value = 'Congo'
$ws.Range("B113").Value = 3546
$ws.Range("C113").Value = 346
$ws.Range("D113").Value = 1589
$ws.Range("E113").Value = 1899
$ws.Range("G113").Value = 4

# Row 114: 'Montenegro' -> 'Tailandia'
$ws.Range("A114").Value = 'Tailandia'
$ws.Range("B114").Value = 3320
$ws.Range("C114").Value = 3
$ws.Range("D114").Value = 3142
$ws.Range("E114").Value = 120
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 58

# Row 115: 'Somalia' -> 'Montenegro'
$ws.Range("A115").Value = 'Montenegro'
$ws.Range("B115").Value = 3301
$ws.Range("C115").Value = 43
$ws.Range("D115").Value = 1445
$ws.Range("E115").Value = 1804
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 52

# Row 116: 'Congo' -> 'Somalia'
$ws.Range("A116").Value = 'Somalia'
$ws.Range("B116").Value = 3220
$ws.Range("D116").Value = 1598
$ws.Range("E116").Value = 1529
$ws.Range("H116").Value = 93

# Row 121: 'Cabo Verde'
$ws.Range("B121").Value = 2583
$ws.Range("C121").Value = 36
$ws.Range("D121").Value = 1911
$ws.Range("E121").Value = 647
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 25

# Row 140: 'Letonia' -> 'Angola'
$ws.Range("A140").Value = 'Angola'
$ws.Range("B140").Value = 1280
$ws.Range("C140").Value = 81
$ws.Range("D140").Value = 476
$ws.Range("E140").Value = 746
$ws.Range("G140").Value = 3
$ws.Range("H140").Value = 58

# Row 141: 'Jordania' -> 'Letonia'
$ws.Range("A141").Value = 'Letonia'
$ws.Range("B141").Value = 1246
$ws.Range("C141").Value = 3
$ws.Range("D141").Value = 1052
$ws.Range("E141").Value = 162
$ws.Range("H141").Value = 32

# Row 142: 'Liberia' -> 'Jordania'
$ws.Range("A142").Value = 'Jordania'
$ws.Range("B142").Value = 1218
$ws.Range("C142").Value = 5
$ws.Range("D142").Value = 1131
$ws.Range("E142").Value = 76
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 11

# Row 143: 'Angola' -> 'Liberia'
$ws.Range("A143").Value = 'Liberia'
$ws.Range("B143").Value = 1214
$ws.Range("C143").Value = 7
$ws.Range("D143").Value = 696
$ws.Range("E143").Value = 440
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 78

# Row 151: 'Republica del Chad'
$ws.Range("D151").Value = 814
$ws.Range("E151").Value = 47

# Row 174: 'Martinica' -> 'Guadalupe'
$ws.Range("A174").Value = 'Guadalupe'
$ws.Range("B174").Value = 272
$ws.Range("C174").Value = 7
$ws.Range("D174").Value = 179
$ws.Range("E174").Value = 79
$ws.Range("H174").Value = 14

# Row 175: 'Guadalupe' -> 'Martinica'
$ws.Range("A175").Value = 'Martinica'
$ws.Range("B175").Value = 269
$ws.Range("D175").Value = 98
$ws.Range("E175").Value = 156
$ws.Range("H175").Value = 15

# Row 180: 'Trinidad yTobago'
$ws.Range("D180").Value = 135
$ws.Range("E180").Value = 39

# Row 185: 'Aruba'
$ws.Range("B185").Value = 124
$ws.Range("C185").Value = 2
$ws.Range("D185").Value = 112
$ws.Range("E185").Value = 9

